$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "333.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.40%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.30%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.695"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.16%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08387"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.86%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.819"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.57%"

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.527"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.20%"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.963"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.04%"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.814"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.48%"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9483"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.81%"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1222"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.65%"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1970"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.64%"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1015"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "8.55%"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04508"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "15.76%"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1068"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.63%"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001286"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.39%"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005895"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-5.45%"

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.476"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.07%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.07%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.757"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.13%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1349"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.78%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2720"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "11.11%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04396"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.13%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001237"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.49%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004332"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.26%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001233"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2.58%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003997"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "31.12%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02915"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.63%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05858"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.57%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007918"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.80%"

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Dexo"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009382"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.74%"

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1425"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.05%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002146"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.16%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009907"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-16.43%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007604"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "9.08%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003183"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.02%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002275"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.25%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
